# Updated cryptos list on Tue Sep 26 01:37:40 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.275.97"
$ws.Range("E2").Value = "  -0.45%  "
$ws.Range("D3").Value = "1.588.57"
$ws.Range("E3").Value = "  -0.05%  "
$ws.Range("E4").Value = "  -0.33%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.26"
$ws.Range("E5").Value = "  +0.19%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.504"
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("E7").Value = "  -0.30%  "
$ws.Range("E8").Value = "  +0.20%  "
$ws.Range("E9").Value = "  -0.61%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.35"
$ws.Range("E10").Value = "  -0.87%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0846"
$ws.Range("E11").Value = "  +0.49%  "
$ws.Range("D12").Value = "1.813.03"
$ws.Range("E12").Value = "  +0.02%  "
$ws.Range("D13").Value = "1.599.60"
$ws.Range("E13").Value = "  +0.88%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.08"
$ws.Range("E14").Value = "  +1.16%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.521"
$ws.Range("E15").Value = "  +0.65%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.45"
$ws.Range("E16").Value = "  +0.14%  "
$ws.Range("D17").Value = "26.291.16"
$ws.Range("E17").Value = "  -0.33%  "
$ws.Range("E18").Value = "  -1.13%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.48"
$ws.Range("E19").Value = "  +3.66%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "210.89"
$ws.Range("E20").Value = "  +1.77%  "
$ws.Range("E21").Value = "  -0.28%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.28"
$ws.Range("E22").Value = "  -0.17%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.98"
$ws.Range("E23").Value = "  +1.79%  "
$ws.Range("E24").Value = "  -3.55%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.05"
$ws.Range("E25").Value = "  -0.45%  "
$ws.Range("E26").Value = "  -0.31%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.06"
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.112"
$ws.Range("E28").Value = "  -1.19%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.21"
$ws.Range("E29").Value = "  -0.82%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0503"
$ws.Range("E30").Value = "  -0.29%  "
$ws.Range("E31").Value = "  +0.49%  "
$ws.Range("E32").Value = "  -1.10%  "
$ws.Range("E33").Value = "  +1.04%  "
$ws.Range("D34").Value = "1.321.47"
$ws.Range("E34").Value = "  +2.69%  "
$ws.Range("E35").Value = "  -1.55%  "
$ws.Range("E36").Value = "  +1.02%  "
$ws.Range("E37").Value = "  -0.67%  "
$ws.Range("E38").Value = "  -0.29%  "
$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.811"
$ws.Range("E39").Value = "  -1.20%  "
$ws.Range("B40").Value = "PaxDollar"
$ws.Range("C40").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  -0.32%  "
$ws.Range("B41").Value = "WEMIXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.02"
$ws.Range("E41").Value = "  -20.30%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.68"
$ws.Range("E42").Value = "  +4.20%  "
$ws.Range("B43").Value = "MXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.14"
$ws.Range("E43").Value = "  -0.04%  "
$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.765"
$ws.Range("E44").Value = "  -0.42%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "61.89"
$ws.Range("E45").Value = "  -1.18%  "
$ws.Range("D46").Value = "1.724.77"
$ws.Range("E46").Value = "  -0.17%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "87.85"
$ws.Range("E47").Value = "  -1.42%  "
$ws.Range("D48").Value = "0.0₆0104"
$ws.Range("E48").Value = "  -1.86%  "
$ws.Range("E49").Value = "  -5.62%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0504"
$ws.Range("E50").Value = "  -1.22%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0974"
$ws.Range("E51").Value = "  -4.88%  "
